$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: replace with old row 12 data (column A / id sequence unchanged) ---
$ws.Range("B11").Value = 6227815
$ws.Range("C11").Value = "Canada Premier League"
$ws.Range("D11").Value = 45094.625
$ws.Range("E11").Value = "HFX Wanderers"
$ws.Range("F11").Value = "Cavalry FC"
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = "H"
$ws.Range("J11").Value = 2.6
$ws.Range("K11").Value = 3.2
$ws.Range("L11").Value = 2.4
$ws.Range("M11").Value = 3.3
$ws.Range("N11").Value = 3
$ws.Range("O11").Value = 2.15
$ws.Range("P11").Value = 0.25
$ws.Range("Q11").Value = 1.925
$ws.Range("R11").Value = 1.875
$ws.Range("S11").Value = 2.25
$ws.Range("T11").Value = 2
$ws.Range("U11").Value = 1.8
$ws.Range("V11").Value = 2.3
$ws.Range("W11").Value = -1
$ws.Range("X11").Value = -1
$ws.Range("Y11").Value = 0.925
$ws.Range("Z11").Value = -1
$ws.Range("AA11").Value = 1
$ws.Range("AB11").Value = -1

# --- Row 12: replace with old row 11 data (column A / id sequence unchanged) ---
$ws.Range("B12").Value = 6240280
$ws.Range("C12").Value = "Canada Premier League"
$ws.Range("D12").Value = 45094.625
$ws.Range("E12").Value = "Atletico Ottawa"
$ws.Range("F12").Value = "Vancouver FC"
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = "H"
$ws.Range("J12").Value = 1.571
$ws.Range("K12").Value = 3.4
$ws.Range("L12").Value = 5.5
$ws.Range("M12").Value = 1.444
$ws.Range("N12").Value = 3.8
$ws.Range("O12").Value = 6
$ws.Range("P12").Value = -1.25
$ws.Range("Q12").Value = 1.95
$ws.Range("R12").Value = 1.85
$ws.Range("S12").Value = 2.75
$ws.Range("T12").Value = 1.975
$ws.Range("U12").Value = 1.825
$ws.Range("V12").Value = 0.444
$ws.Range("W12").Value = -1
$ws.Range("X12").Value = -1
$ws.Range("Y12").Value = -0.5
$ws.Range("Z12").Value = 0.425
$ws.Range("AA12").Value = -1
$ws.Range("AB12").Value = 0.825

# --- Row 83: replace with old row 84 data (column A / id sequence unchanged) ---
$ws.Range("B83").Value = 7301364
$ws.Range("C83").Value = "Canada Premier League"
$ws.Range("D83").Value = 45206.75
$ws.Range("E83").Value = "Forge FC"
$ws.Range("F83").Value = "Atletico Ottawa"
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 1
$ws.Range("I83").Value = "A"
$ws.Range("J83").Value = 1.8
$ws.Range("K83").Value = 3.6
$ws.Range("L83").Value = 3.5
$ws.Range("M83").Value = 1.533
$ws.Range("N83").Value = 3.8
$ws.Range("O83").Value = 5
$ws.Range("P83").Value = -1
$ws.Range("Q83").Value = 1.975
$ws.Range("R83").Value = 1.825
$ws.Range("S83").Value = 2.5
$ws.Range("T83").Value = 1.9
$ws.Range("U83").Value = 1.9
$ws.Range("V83").Value = -1
$ws.Range("W83").Value = -1
$ws.Range("X83").Value = 4
$ws.Range("Y83").Value = -1
$ws.Range("Z83").Value = 0.825
$ws.Range("AA83").Value = -1
$ws.Range("AB83").Value = 0.8999999999999999

# --- Row 84: replace with old row 83 data (column A / id sequence unchanged) ---
$ws.Range("B84").Value = 6227884
$ws.Range("C84").Value = "Canada Premier League"
$ws.Range("D84").Value = 45206.75
$ws.Range("E84").Value = "Cavalry FC"
$ws.Range("F84").Value = "Pacific FC CA"
$ws.Range("G84").Value = 3
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = "H"
$ws.Range("J84").Value = 2.25
$ws.Range("K84").Value = 3.1
$ws.Range("L84").Value = 2.875
$ws.Range("M84").Value = 2.05
$ws.Range("N84").Value = 3.2
$ws.Range("O84").Value = 3.2
$ws.Range("P84").Value = -0.25
$ws.Range("Q84").Value = 1.825
$ws.Range("R84").Value = 1.975
$ws.Range("S84").Value = 2.5
$ws.Range("T84").Value = 1.825
$ws.Range("U84").Value = 1.975
$ws.Range("V84").Value = 1.05
$ws.Range("W84").Value = -1
$ws.Range("X84").Value = -1
$ws.Range("Y84").Value = 0.825
$ws.Range("Z84").Value = -1
$ws.Range("AA84").Value = 0.825
$ws.Range("AB84").Value = -1

# --- Row 108: rewrite with new match data (A/D already have correct styles) ---
$ws.Range("B108").Value = 7802943
$ws.Range("C108").Value = "Canada Premier League"
$ws.Range("D108").Value = 45429.95833333334
$ws.Range("E108").Value = "Pacific FC CA"
$ws.Range("F108").Value = "Atletico Ottawa"
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 1
$ws.Range("I108").Value = "A"
$ws.Range("J108").Value = 2
$ws.Range("K108").Value = 3.4
$ws.Range("L108").Value = 3.2
$ws.Range("M108").Value = 2.25
$ws.Range("N108").Value = 3.25
$ws.Range("O108").Value = 2.8
$ws.Range("P108").Value = -0.25
$ws.Range("Q108").Value = 2.025
$ws.Range("R108").Value = 1.775
$ws.Range("S108").Value = 2.25
$ws.Range("T108").Value = 1.975
$ws.Range("U108").Value = 1.825
$ws.Range("V108").Value = -1
$ws.Range("W108").Value = -1
$ws.Range("X108").Value = 1.8
$ws.Range("Y108").Value = -1
$ws.Range("Z108").Value = 0.7749999999999999
$ws.Range("AA108").Value = -1
$ws.Range("AB108").Value = 0.825

# --- New row 109 ---
$ws.Range("A109").Value = 107
$ws.Range("B109").Value = 7802879
$ws.Range("C109").Value = "Canada Premier League"
$ws.Range("D109").Value = 45430.70833333334
$ws.Range("E109").Value = "Forge FC"
$ws.Range("F109").Value = "Vancouver FC"
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 2
$ws.Range("I109").Value = "A"
$ws.Range("J109").Value = 1.6
$ws.Range("K109").Value = 4
$ws.Range("L109").Value = 4.333
$ws.Range("M109").Value = 1.55
$ws.Range("N109").Value = 4
$ws.Range("O109").Value = 4.5
$ws.Range("P109").Value = -1
$ws.Range("Q109").Value = 1.975
$ws.Range("R109").Value = 1.825
$ws.Range("S109").Value = 2.75
$ws.Range("T109").Value = 1.975
$ws.Range("U109").Value = 1.825
$ws.Range("V109").Value = -1
$ws.Range("W109").Value = -1
$ws.Range("X109").Value = 3.5
$ws.Range("Y109").Value = -1
$ws.Range("Z109").Value = 0.825
$ws.Range("AA109").Value = 0.4875
$ws.Range("AB109").Value = -0.5

# --- New row 110 ---
$ws.Range("A110").Value = 108
$ws.Range("B110").Value = 7803367
$ws.Range("C110").Value = "Canada Premier League"
$ws.Range("D110").Value = 45430.83333333334
$ws.Range("E110").Value = "Cavalry FC"
$ws.Range("F110").Value = "York United FC"
$ws.Range("G110").Value = 2
$ws.Range("H110").Value = 2
$ws.Range("I110").Value = "D"
$ws.Range("J110").Value = 2
$ws.Range("K110").Value = 3.4
$ws.Range("L110").Value = 3.2
$ws.Range("M110").Value = 1.65
$ws.Range("N110").Value = 3.7
$ws.Range("O110").Value = 4.333
$ws.Range("P110").Value = -0.75
$ws.Range("Q110").Value = 1.85
$ws.Range("R110").Value = 1.95
$ws.Range("S110").Value = 2.5
$ws.Range("T110").Value = 1.825
$ws.Range("U110").Value = 1.975
$ws.Range("V110").Value = -1
$ws.Range("W110").Value = 2.7
$ws.Range("X110").Value = -1
$ws.Range("Y110").Value = -1
$ws.Range("Z110").Value = 0.95
$ws.Range("AA110").Value = 0.825
$ws.Range("AB110").Value = -1

# --- New row 111 ---
$ws.Range("A111").Value = 109
$ws.Range("B111").Value = 7803368
$ws.Range("C111").Value = "Canada Premier League"
$ws.Range("D111").Value = 45432.66666666666
$ws.Range("E111").Value = "HFX Wanderers"
$ws.Range("F111").Value = "Valour FC"
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 2
$ws.Range("I111").Value = "A"
$ws.Range("J111").Value = 2.1
$ws.Range("K111").Value = 3.2
$ws.Range("L111").Value = 3.2
$ws.Range("M111").Value = 1.7
$ws.Range("N111").Value = 3.6
$ws.Range("O111").Value = 4.1
$ws.Range("P111").Value = -0.75
$ws.Range("Q111").Value = 1.95
$ws.Range("R111").Value = 1.85
$ws.Range("S111").Value = 2.5
$ws.Range("T111").Value = 1.95
$ws.Range("U111").Value = 1.85
$ws.Range("V111").Value = -1
$ws.Range("W111").Value = -1
$ws.Range("X111").Value = 3.1
$ws.Range("Y111").Value = -1
$ws.Range("Z111").Value = 0.8500000000000001
$ws.Range("AA111").Value = 0.95
$ws.Range("AB111").Value = -1

# --- Copy formatting (border/bold/centered for column A, date format for column D) for the new rows from row 108 ---
$ws.Range("A108").Copy()
$ws.Range("A109:A111").PasteSpecial(-4122)
$ws.Range("D108").Copy()
$ws.Range("D109:D111").PasteSpecial(-4122)
$excel.CutCopyMode = $false
